$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("114_1")

$ws.Range("B5").Value = 42369
$ws.Range("C5").Value = 42369
$ws.Range("D5").Value = 42369
